# The deck had slides 5, 7, 9 and 10 marked as hidden (<p:sld show="0">).
# Un-hide them so they play normally again during the slide show.
$p = $ppt.ActivePresentation

$slideNumbers = @(5, 7, 9, 10)

foreach ($n in $slideNumbers) {
    $slide = $p.Slides.Item($n)
    $slide.SlideShowTransition.Hidden = $false
}
